$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (the "moon / first animal" record) -- update the backward-reasoning
# related columns (bck_final_ans, bck_final_ans_exp, bck_final_question,
# original_response)
$ws.Range("F2").Value = "The conclusion is not plausible."
$ws.Range("G2").Value = "The conclusion that Neil Armstrong was the first animal to land on the Moon is not plausible because humans are not typically referred to as animals in the context of space exploration. Additionally, historical records and scientific knowledge confirm that no animals were sent to the Moon, contradicting the conclusion drawn from the observation. Therefore, the conclusion is not logically and factually coherent with our understanding of reality."
$ws.Range("H2").Value = "Was Neil Armstrong the first animal to land on the moon?"
$ws.Range("I2").Value = "{'Answer:': 'No animals were ever sent to the Moon. Although, since humans are technically animals, one could say that the first animal sent to the Moon was Neil Armstrong. He belonged to the species Homo sapiens.', 'Source:': 'Study.com (https://homework.study.com/explanation/what-was-the-first-animal-to-land-on-the-moon.html#:~:text=Answer%20and%20Explanation%3A,to%20the%20species%20Homo%20sapiens.)', 'Premise of the Question:': 'Valid', 'Explanation:': ''}"

# Row 3 (the "Leonardo DiCaprio" record) -- update fwd_final_ans plus the
# backward-reasoning related columns
$ws.Range("D3").Value = "Leonardo DiCaprio does not have any children, so he does not have a third child."
$ws.Range("F3").Value = "The conclusion is highly plausible as Leonardo DiCaprio does not have any children."
$ws.Range("G3").Value = "Leonardo DiCaprio has never been married and has no children, as confirmed by reliable sources. Additionally, there is no credible information suggesting that he has any children. Therefore, the conclusion that he does not have a third child is highly credible and aligns with the available evidence and general knowledge."
$ws.Range("H3").Value = "Does Leonardo DiCaprio have any children?"
$ws.Range("I3").Value = "{'Answer:': 'Leonardo DiCaprio does not have any children, so he does not have a third child.', 'Source:': 'Quora (https://www.quora.com/How-many-women-has-Leonardo-DiCaprio-dated-How-many-did-he-marry-and-how-many-children-does-he-have)', 'Premise of the Question:': 'Invalid', 'Explanation:': ''}"

$wb.Save()
